$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to text
# so Excel does not auto-convert them to numeric values (losing the original
# text formatting, e.g. trailing zeros).
$textForceRefs = @("D5", "D6", "D7", "D10", "D11", "D15", "D16", "D21", "D23", "D26", "D30", "D31", "D33", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply all cell value updates
$ws.Range("D2").Value = "41.509.45"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.203.04"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "239.79"
$ws.Range("D6").Value = "0.618"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").Value = "71.20"
$ws.Range("E7").Value = "  -3.67%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -4.58%  "
$ws.Range("D10").Value = "40.98"
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").Value = "0.0934"
$ws.Range("E11").Value = "  -4.12%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  -4.46%  "
$ws.Range("D14").Value = "2.539.88"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "14.03"
$ws.Range("E15").Value = "  -2.64%  "
$ws.Range("D16").Value = "0.819"
$ws.Range("E16").Value = "  -3.80%  "
$ws.Range("D17").Value = "2.207.66"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").Value = "41.486.71"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("E19").Value = "  -11.40%  "
$ws.Range("E20").Value = "  -1.43%  "
$ws.Range("D21").Value = "71.10"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("E22").Value = "  +7.16%  "
$ws.Range("D23").Value = "227.11"
$ws.Range("E23").Value = "  -1.75%  "
$ws.Range("E24").Value = "  -6.47%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "11.19"
$ws.Range("E26").Value = "  -5.57%  "
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").Value = "166.14"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").Value = "20.21"
$ws.Range("E31").Value = "  -4.07%  "
$ws.Range("E32").Value = "  -2.32%  "
$ws.Range("D33").Value = "30.13"
$ws.Range("E33").Value = "  +2.64%  "
$ws.Range("E34").Value = "  -7.35%  "
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("E36").Value = "  -9.78%  "
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("D38").Value = "0.0296"
$ws.Range("E38").Value = "  -2.32%  "
$ws.Range("D39").Value = "12.78"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("D41").Value = "5.53"
$ws.Range("E41").Value = "  -1.46%  "
$ws.Range("D42").Value = "62.78"
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "0.192"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("D44").Value = "8.52"
$ws.Range("E44").Value = "  -3.27%  "
$ws.Range("D45").Value = "0.0989"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("D46").Value = "100.48"
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("B47").Value = "HuobiToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D47").Value = "2.88"
$ws.Range("E47").Value = "  +6.05%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("B49").Value = "TrustWalletToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D49").Value = "1.15"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "2.28"
$ws.Range("E50").Value = "  -3.64%  "
$ws.Range("D51").Value = "2.415.34"
$ws.Range("E51").Value = "  -0.90%  "

# Restore default style on the cells we forced to text, keeping them as text
# values but without leaving a stray number-format style behind.
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
